# Apply BOM updates for the I2C0 ESD circuit parts added during design review:
#  - RCLAMP0504FB (U1) row: add U18 as a second instance (ESD protection IC)
#  - 0.1uF capacitor row: add C41
#  - 100 ohm resistor row: add R73, R74, R75

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: 100 ohm resistors, add R73, R74, R75 ; Qty 3 -> 6
$ws.Range("A46").Value = 6
$ws.Range("E46").Value = "R41, R48, R71, R73, R74, R75"

# Row 36: 0.1uF capacitors, add C41 ; Qty 14 -> 15
$ws.Range("A36").Value = 15
$ws.Range("E36").Value = "C4, C11, C16, C17, C18, C20, C21, C24, C25, C26, C28, C31, C39, C40, C41"

# Row 24: RCLAMP0504FB / U1 -> U1, U18 ; Qty 1 -> 2
$ws.Range("A24").Value = 2
$ws.Range("E24").Value = "U1, U18"
